# Updated symbol list (price / 1h volume change) refresh.
# Values are numeric- / percent-looking text that must stay as literal text
# (matching the source data's inline-string cells), so each value is written
# with a leading apostrophe (forces Excel to treat it as text instead of
# auto-converting to a Number/Percentage) and then ClearFormats() is used to
# drop the "quote prefix" cell style that the apostrophe trick leaves behind,
# so the cell keeps the workbook's original (default/General) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "289.85";  "E2"  = "-4.44%"
    "D3"  = "30.75";   "E3"  = "-4.46%"
    "D4"  = "4.947";   "E4"  = "0.51%"
    "D5"  = "0.07161"; "E5"  = "-8.95%"
    "D6"  = "1.794";   "E6"  = "-12.70%"
    "D7"  = "7.657";   "E7"  = "-2.30%"
    "D8"  = "3.738";   "E8"  = "-2.82%"
    "D9"  = "0.8966";  "E9"  = "-3.07%"
    "D10" = "0.1647";  "E10" = "-6.67%"
    "D11" = "0.07675"; "E11" = "-1.80%"
    "D12" = "0.08046"; "E12" = "-6.23%"
    "D13" = "0.03031"; "E13" = "-4.07%"
                        "E14" = "-0.31%"
    "D15" = "0.001507"; "E15" = "-1.24%"
    "D16" = "0.005767"; "E16" = "-2.00%"
    "D18" = "3.475";    "E18" = "0.30%"
                         "E19" = "-1.64%"
                         "E20" = "0.00%"
    "D21" = "0.1330";   "E21" = "1.05%"
    "D22" = "4.045";    "E22" = "-5.76%"
    "D23" = "0.1999";   "E23" = "0.42%"
    "D24" = "0.04516";  "E24" = "-1.30%"
                         "E25" = "-0.87%"
    "D26" = "0.003998"; "E26" = "-10.31%"
    "D27" = "0.0001250"; "E27" = "0.00%"
    "D39" = "0.01608";  "E39" = "-7.36%"
    "D40" = "0.04369";  "E40" = "-9.07%"
    "D41" = "0.007380"; "E41" = "-1.81%"
    "D42" = "0.1305";   "E42" = "-4.41%"
    "D43" = "0.002061"; "E43" = "-12.71%"
    "D44" = "0.009534"; "E44" = "-9.85%"
    "D45" = "0.00005978"; "E45" = "-2.93%"
                           "E46" = "-0.04%"
    "D47" = "2.246";    "E47" = "172.73%"
    "D48" = "0.003000"; "E48" = "-3.17%"
    "D49" = "0.00002100"; "E49" = "-0.04%"
    "D50" = "0.0002000";  "E50" = "-0.04%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.ClearFormats()
}
